# Commit: "Remove unused code and update template"
#
# tbl_spesifikasi ("sheet4.xml") had a single "os" column (E) between
# "processor" and "memory". The template is updated so that single column
# is replaced by three columns "os1", "os2", "os3". This also shifts the
# active/selected sheet from tbl_aset to tbl_spesifikasi.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl_spesifikasi")

# Insert two new columns right after the existing "os" column (E) so the
# old column becomes the first of three, shifting memory/hard_drive/
# keterangan (and the legend column further right) two columns over -
# matching the dimension/merge/col-width changes in the diff.
$ws.Range("F1:G1").EntireColumn.Insert() | Out-Null

# Relabel the (now) three os-related headers.
$ws.Range("E2").Value = "os1"
$ws.Range("F2").Value = "os2"
$ws.Range("G2").Value = "os3"

# Match the new selection on this sheet and make it the active tab
# (previously tbl_aset was the active/selected sheet).
$ws.Range("G2").Select() | Out-Null
$ws.Activate() | Out-Null
